$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 502
$ws.Range("F4").Value = 434
$ws.Range("F5").Value = 8625
$ws.Range("F6").Value = 13
$ws.Range("F7").Value = 10856
$ws.Range("F18").Value = 81
$ws.Range("F20").Value = 415
$ws.Range("F24").Value = 582
$ws.Range("F27").Value = 70
$ws.Range("F30").Value = 1221
$ws.Range("F38").Value = 291
$ws.Range("F41").Value = 522
$ws.Range("F42").Value = 357
$ws.Range("F45").Value = 645
$ws.Range("F48").Value = 117

$ws = $wb.Worksheets.Item(2)
$ws.Range("F12").Value = 64
$ws.Range("F17").Value = 388

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 213
$ws.Range("F3").Value = 2818

$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 213
$ws.Range("F4").Value = 502
$ws.Range("F8").Value = 434
$ws.Range("F9").Value = 8625
$ws.Range("F10").Value = 13
$ws.Range("F11").Value = 10856
$ws.Range("F16").Value = 81
$ws.Range("F18").Value = 415
$ws.Range("F21").Value = 582
$ws.Range("F24").Value = 70
$ws.Range("F28").Value = 1221
$ws.Range("F36").Value = 64
$ws.Range("F39").Value = 522
$ws.Range("F40").Value = 357
$ws.Range("F45").Value = 388
$ws.Range("F46").Value = 645
$ws.Range("F49").Value = 117
